$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from an existing header cell (H1) so the new headers match
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data for columns I (I0) and J (IF), rows 2-33
$data = @(
    @(2, 7),
    @(4, 5),
    @(5, 5),
    @(5, 6),
    @(1, 5),
    @(8, 8),
    @(3, 7),
    @(1, 5),
    @(6, 7),
    @(8, 8),
    @(2, 5),
    @(8, 8),
    @(8, 9),
    @(7, 9),
    @(1, 7),
    @(1, 4),
    @(5, 7),
    @(1, 4),
    @(5, 7),
    @(8, 8),
    @(1, 5),
    @(1, 5),
    @(2, 5),
    @(1, 4),
    @(1, 4),
    @(2, 5),
    @(6, 9),
    @(1, 3),
    @(5, 6),
    @(8, 8),
    @(7, 9),
    @(4, 5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
